# Model and template update: add a third header row (row 3) to the
# "data_info(zymography)" sheet holding the enum/description tags that
# correspond to each existing column header in row 1:
#   Operator                          -> #Manipulateur
#   SampleID                          -> #Desc:IdentifiantEchantillon
#   Date                              -> #Date
#   LaboratoryOperatingMode           -> #ModeOderatoireLaboratoire
#   CriticalApparatusCriticalSoftware -> #AppareilLogicielCritique
#   CriticalProduct                   -> #ProduitCritique
#   RawDataPathway                    -> #LieuStockageDonneesBrutes
#   (the remaining columns H:O get an empty string marker)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"

# H3:O3 hold an explicit empty-text marker (not a blank/cleared cell).
# A bare "" assignment clears the cell instead of keeping it as a typed
# empty string, so enter it as a forced-text (leading apostrophe) value
# and then strip the resulting "quote prefix" number format so the cell
# keeps its default (unstyled) appearance, same as every other cell here.
$ws.Range("H3:O3").Value = "'"
$ws.Range("H3:O3").ClearFormats()
